$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 currently holds Luis Armando Cordoba Carmona's data; row 17 holds
# Jessica Rodriguez Pereira's data. Swap the two rows' content so that
# Jessica's record now appears first (row 16) and Luis's second (row 17).
$ws.Range("C16").Value = "45559364"
$ws.Range("D16").Value = "JESSICA RODRIGUEZ PEREIRA"
$ws.Range("E16").Value = "2211"
$ws.Range("F16").Value = 20000
$ws.Range("G16").Value = 1000000

$ws.Range("C17").Value = "1047424394"
$ws.Range("D17").Value = "LUIS ARMANDO CORDOBA CARMONA"
$ws.Range("E17").Value = "2302"
$ws.Range("F17").Value = 1547
$ws.Range("G17").Value = 1160000
